$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Documentacion - Informe Final SQA
# Insert 8 new rows starting at row 11 to make room for the new nomenclature
# entries related to Quality Management (SQA). This shifts the existing
# E201..E211 rows down from rows 11-21 to rows 19-29.
$ws.Rows("11:18").Insert()

# Row 10: the "Plan de Pruebas" document is replaced here by the new
# "Manejo del Ambiente Controlado" document.
$ws.Range("B10").Value = "Manejo del Ambiente Controlado"

# Row 11: new entry for the Informe Final de SQA document.
$ws.Range("A11").Value = "E1010"
$ws.Range("B11").Value = "Informe Final de SQA"

# Row 12: new code only, no document name yet.
$ws.Range("A12").Value = "E1011"

# Row 13: new code; "Plan de Pruebas" is relocated here.
$ws.Range("A13").Value = "E1012"
$ws.Range("B13").Value = "Plan de Pruebas"

# Rows 14-18: remaining new codes reserved for future documents (column A only).
$ws.Range("A14").Value = "E1013"
$ws.Range("A15").Value = "E1014"
$ws.Range("A16").Value = "E1015"
$ws.Range("A17").Value = "E1016"
$ws.Range("A18").Value = "E1017"

# Update the active cell selection to match the edited workbook state.
$ws.Range("C16").Select() | Out-Null
